# Connected projected to PM and updated test plans. Added jenkins build file
#
# 1) On the existing "OppData" sheet, add a new "Probability" column (D)
#    with values for the first couple of opportunity rows, and move the
#    active selection to D2.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("OppData")

$ws1.Range("D1").Value = "Probability"
$ws1.Range("D2").Value = 10
$ws1.Range("D3").Value = 15

$ws1.Columns.Item(4).ColumnWidth = 26.45

# 2) Add a new worksheet ("Sheet1") after OppData with a Rule/Conversion
#    lookup table, and leave it as the active sheet/tab.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"

$ws2.Range("A1").Value = "Rule"
$ws2.Range("B1").Value = "Conversion"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "list_amount > 50000"

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = 'business_unit = "Transportation"'

$ws2.Range("A4").Value = 3
$ws2.Range("A5").Value = 4
$ws2.Range("A6").Value = 5
$ws2.Range("A7").Value = 6
$ws2.Range("A8").Value = 7

$ws2.Columns.Item(2).ColumnWidth = 29.18

# Restore the selection on OppData to D2 (Add()/typing above moved it),
# then re-select B3 on the new sheet so it ends up the active tab/cell.
$ws1.Range("D2").Select() | Out-Null
$ws2.Range("B3").Select() | Out-Null
